# Board inventory, ordered some boards and stencils
# Change every "PEND" (pending) status on the Boards sheet to "ORDR" (ordered),
# applying a distinct yellow fill (the green fill is used for "RCVD").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Boards")

$cellsToUpdate = @("E6", "E7", "E8", "B9", "E9", "B11", "E11", "B12", "E12", "B13", "E13", "B14", "E14")

foreach ($addr in $cellsToUpdate) {
    $cell = $ws.Range($addr)
    $cell.Value = "ORDR"
    $cell.Interior.Color = 65535
    $cell.HorizontalAlignment = -4108
}

# Switch the active tab from Inventory to Boards, and move the selection.
$ws.Activate() | Out-Null
$ws.Range("A20").Select() | Out-Null
